# Apply the "clientes" sheet restructuring described by the commit:
# turn row 1 into a new header row (ESTADO/RUT/CLIENTE/DIRECCION/COMUNA/
# TELEFONO/GPS/OTRO), keep row 2 as the normalized client record, and
# add row 3 with the original (pre-normalization) client record plus a
# new "OTRO" test marker in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clientes")

# --- Row 1: new header labels (column H is new) ---
$ws.Range("A1").Value = "ESTADO"
$ws.Range("B1").Value = "RUT"
$ws.Range("C1").Value = "CLIENTE"
$ws.Range("D1").Value = "DIRECCION"
$ws.Range("E1").Value = "COMUNA"
$ws.Range("F1").Value = "TELEFONO"
$ws.Range("G1").Value = "GPS"
$ws.Range("H1").Value = "OTRO"

# --- Row 2: re-assert the existing normalized client record values
# (phone number must stay text, not become a number). ---
$ws.Range("A2").Value = "16.742.249-7"
$ws.Range("B2").Value = "Isaias Beroiza Mora"
$ws.Range("C2").Value = "colaco sn km3 parcela 9"
$ws.Range("D2").Value = "Calbuco"
$ws.Range("E2").Value = "'56988809704"
$ws.Range("F2").Value = "por buscar"
$ws.Range("G2").Value = "ok"

# --- Row 3: brand-new row holding the original (pre-normalization)
# client record plus a new "OTRO" test marker in column H. ---
$ws.Range("A3").Value = "activo"
$ws.Range("B3").Value = "16.742.249-7"
$ws.Range("C3").Value = "Isaias Andres Beroiza Mora"
$ws.Range("D3").Value = "colaco sn km3 parcela 9"
$ws.Range("E3").Value = "Calbuco"
$ws.Range("F3").Value = "'56988809704"
$ws.Range("G3").Value = "buscando coordenadas"
$ws.Range("H3").Value = "Cliente prueba"
